# added more games, sped up simulate game logic, and drafted optimization logic
# -> update the team-specific transition-probability matrix on Sheet1 with the
#    freshly simulated (non-zero) probabilities for the affected state cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.5
$ws.Range("P2").Value = 0.2
$ws.Range("S2").Value = 0.1
$ws.Range("J3").Value = 0.2
$ws.Range("P3").Value = 0.6
$ws.Range("S3").Value = 0.2
$ws.Range("F6").Value = 0.1
$ws.Range("J6").Value = 0.2
$ws.Range("Q6").Value = 0.1
$ws.Range("S6").Value = 0.6
$ws.Range("F7").Value = 0.4
$ws.Range("Q7").Value = 0.4
$ws.Range("S7").Value = 0.2
$ws.Range("B8").Value = 0.1
$ws.Range("F8").Value = 0.05
$ws.Range("Q8").Value = 0.25
$ws.Range("R8").Value = 0.15
$ws.Range("S8").Value = 0.45
$ws.Range("B9").Value = 0.25
$ws.Range("J9").Value = 0.5
$ws.Range("Q9").Value = 0.25
$ws.Range("B10").Value = 0.09259259259259259
$ws.Range("F10").Value = 0.05555555555555555
$ws.Range("J10").Value = 0.09259259259259259
$ws.Range("Q10").Value = 0.3518518518518519
$ws.Range("R10").Value = 0.03703703703703703
$ws.Range("S10").Value = 0.3703703703703703
$ws.Range("G11").Value = 0.25
$ws.Range("K11").Value = 0.4166666666666667
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.5
$ws.Range("J15").Value = 0.5
$ws.Range("K15").Value = 0.3333333333333333
$ws.Range("S15").Value = 0.1666666666666667
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.2
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.2
$ws.Range("F17").Value = 0.03703703703703703
$ws.Range("H17").Value = 0.1481481481481481
$ws.Range("I17").Value = 0.03703703703703703
$ws.Range("J17").Value = 0.4814814814814815
$ws.Range("K17").Value = 0.07407407407407407
$ws.Range("O17").Value = 0.07407407407407407
$ws.Range("S17").Value = 0.1481481481481481
$ws.Range("F18").Value = 0.2
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.2
$ws.Range("S18").Value = 0.2
$ws.Range("H19").Value = 0.3
$ws.Range("I19").Value = 0.04
$ws.Range("J19").Value = 0.44
$ws.Range("K19").Value = 0.04
$ws.Range("O19").Value = 0.06
$ws.Range("S19").Value = 0.12
